# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The "K" column (column G) values are recalculated/regenerated; write the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 2
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 3
    14 = 2
    15 = 1
    16 = 1
    17 = 0
    18 = 0
    21 = 2
    22 = 1
    24 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
